$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Pischke & von Wachter (2008) - parents' effect
$ws.Range("K3").Value = "Pischke & von Wachter (2008), Table II, Microzensus Basic Track IV"
$ws.Range("A3").Value = "lwage_effect_parents_pw"
$ws.Range("B3").Value = 0.013
$ws.Range("C3").Value = 0.011
$ws.Range("J3").Value = 2

# Row 4: Kamhöfer & Schmitz (2016) - parents' effect
$ws.Range("A4").Value = "lwage_effect_parents_ks"
$ws.Range("B4").Value = -0.0004
$ws.Range("C4").Value = 0.0276
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = "Kamhöfer & Schmitz (2016), Table I, Basic"

# Update selection to A3 as in the diff
$ws.Range("A3").Select()
